# Updates the crypto price/volume table to the latest scraped snapshot
# (GitHub Actions data refresh), including two coin re-rankings that swapped
# row order (MultiversX/LidoDAOToken at 41/42, Stacks/WOONetwork at 50/51).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some new Price values are plain numerics (e.g. "175.50"); the source sheet
# stores every Price/Volume cell as text, so force a text format first on the
# cells that would otherwise be auto-coerced into Excel numbers (which would
# silently drop things like trailing zeros or significant digits).
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D51').NumberFormat = '@'

$ws.Range('D2').Value = '44.036.72'
$ws.Range('E2').Value = '  +2.60%  '
$ws.Range('D3').Value = '2.247.11'
$ws.Range('E3').Value = '  +1.32%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').Value = '268.91'
$ws.Range('E5').Value = '  +4.60%  '
$ws.Range('D6').Value = '87.72'
$ws.Range('E6').Value = '  +12.90%  '
$ws.Range('E7').Value = '  +1.20%  '
$ws.Range('E8').Value = '  -0.07%  '
$ws.Range('D9').Value = '0.616'
$ws.Range('E9').Value = '  +3.37%  '
$ws.Range('D10').Value = '46.09'
$ws.Range('E10').Value = '  +7.11%  '
$ws.Range('D11').Value = '0.0931'
$ws.Range('E11').Value = '  +2.76%  '
$ws.Range('D12').Value = '7.56'
$ws.Range('E12').Value = '  +7.99%  '
$ws.Range('E13').Value = '  +2.26%  '
$ws.Range('D14').Value = '2.580.26'
$ws.Range('E14').Value = '  +1.12%  '
$ws.Range('D15').Value = '14.98'
$ws.Range('E15').Value = '  +3.84%  '
$ws.Range('D16').Value = '2.242.50'
$ws.Range('E16').Value = '  +1.07%  '
$ws.Range('D17').Value = '0.799'
$ws.Range('E17').Value = '  +1.78%  '
$ws.Range('D18').Value = '44.010.41'
$ws.Range('E18').Value = '  +2.66%  '
$ws.Range('E19').Value = '  +0.57%  '
$ws.Range('E20').Value = '  +1.01%  '
$ws.Range('D21').Value = '70.28'
$ws.Range('E21').Value = '  -1.15%  '
$ws.Range('D22').Value = '2.41'
$ws.Range('E22').Value = '  +5.25%  '
$ws.Range('D23').Value = '233.69'
$ws.Range('E23').Value = '  +1.62%  '
$ws.Range('D24').Value = '8.93'
$ws.Range('E24').Value = '  -4.21%  '
$ws.Range('D25').Value = '2.58'
$ws.Range('E25').Value = '  +17.37%  '
$ws.Range('E26').Value = '  +0.02%  '
$ws.Range('D27').Value = '11.02'
$ws.Range('E27').Value = '  +2.48%  '
$ws.Range('D28').Value = '3.56'
$ws.Range('E28').Value = '  +6.43%  '
$ws.Range('D29').Value = '40.75'
$ws.Range('E29').Value = '  -4.89%  '
$ws.Range('D30').Value = '2.26'
$ws.Range('E30').Value = '  +1.55%  '
$ws.Range('D31').Value = '175.50'
$ws.Range('E31').Value = '  +0.81%  '
$ws.Range('D32').Value = '0.0913'
$ws.Range('E32').Value = '  +4.57%  '
$ws.Range('D33').Value = '20.80'
$ws.Range('E33').Value = '  +1.82%  '
$ws.Range('D34').Value = '5.43'
$ws.Range('E34').Value = '  +3.86%  '
$ws.Range('E35').Value = '  +2.02%  '
$ws.Range('D36').Value = '0.111'
$ws.Range('E36').Value = '  +3.70%  '
$ws.Range('D37').Value = '0.0358'
$ws.Range('E37').Value = '  +0.37%  '
$ws.Range('E38').Value = '  +0.29%  '
$ws.Range('D39').Value = '3.34'
$ws.Range('E39').Value = '  +18.02%  '
$ws.Range('D40').Value = '12.70'
$ws.Range('E40').Value = '  -2.64%  '
$ws.Range('B41').Value = 'LidoDAOToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D41').Value = '2.16'
$ws.Range('E41').Value = '  +1.94%  '
$ws.Range('B42').Value = 'MultiversX'
$ws.Range('C42').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D42').Value = '65.61'
$ws.Range('E42').Value = '  +6.91%  '
$ws.Range('D43').Value = '0.205'
$ws.Range('E43').Value = '  +1.35%  '
$ws.Range('E44').Value = '  +1.50%  '
$ws.Range('E45').Value = '  +2.49%  '
$ws.Range('D46').Value = '8.38'
$ws.Range('E46').Value = '  -0.75%  '
$ws.Range('D47').Value = '100.39'
$ws.Range('E47').Value = '  -3.18%  '
$ws.Range('D48').Value = '1.22'
$ws.Range('E48').Value = '  +7.55%  '
$ws.Range('E49').Value = '  +2.20%  '
$ws.Range('B50').Value = 'WOONetwork'
$ws.Range('C50').Value = 'https://coinranking.com/coin/k-J3YwacF+woonetwork-woo'
$ws.Range('D50').Value = '0.443'
$ws.Range('E50').Value = '  -9.34%  '
$ws.Range('B51').Value = 'Stacks'
$ws.Range('C51').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D51').Value = '1.53'
$ws.Range('E51').Value = '  +4.78%  '
